$wb = $excel.ActiveWorkbook

# 1. Rename "register " sheet to "test_register"
$wsRegister = $wb.Worksheets.Item("register ")
$wsRegister.Name = "test_register"

# 2. Reorder sheet tabs: move "login" to be before "busniess flow"
#    (final order: test_register, login, busniess flow)
$wsLogin = $wb.Worksheets.Item("login")
$wsBusiness = $wb.Worksheets.Item("busniess flow")
$wsLogin.Move($wsBusiness)

# 3. Update the registration-flow test case text in B3:
#    previously "Send and verify registration SMS verification code"
#    now clarified as "press next step to verify registration SMS verification code"
$wsRegister.Range("B3").Value = "press next step to verify registration SMS verification code"

# 4. Reflect the updated view/selection on the test_register sheet
$wsRegister.Activate()
$wsRegister.Range("E12").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
